# Actualización automática 2025-09-24 08:30:08
#
# Updates the three linked sheets (VENTAS POR GRUPO / VENTA MENSUAL /
# CUMPLIMIENTO MENSUAL) of the "GUERRERO FAREZ FABIAN MAURICIO" advisor
# workbook for a handful of revised sale figures.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" --------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# CERAMIKASA S.A.S. (row 14): PIEDRA SINTERIZADA (L) and SAL SOLUBLE (O)
$wsGrupo.Range("L14").Value = 1631.49
$wsGrupo.Range("O14").Value = 2068.07

# FEIJOO MARIN MAURICIO ENRIQUE (row 21): LAVABOS (I)
$wsGrupo.Range("I21").Value = 311.4

# ORTEGA ROMAN LUIS FERNANDO (row 34): PORCELANATO (M)
$wsGrupo.Range("M34").Value = 2948.13

# Footer counter row for SAL SOLUBLE went from 2 non-zero clients to 3
$wsGrupo.Range("O54").Value = "3 de 52"

# --- Sheet "VENTA MENSUAL" -------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# septiembre (F) totals per client, reflecting the VENTAS POR GRUPO changes
$wsMensual.Range("F14").Value = 4558.11
$wsMensual.Range("F21").Value = 5782.28
$wsMensual.Range("F34").Value = 4022.81

# septiembre (F) grand total
$wsMensual.Range("F58").Value = 58377.59

# --- Sheet "CUMPLIMIENTO MENSUAL" -----------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# LAVABOS (row 7)
$wsCumpl.Range("D7").Value = 2349.9
$wsCumpl.Range("E7").Value = -1463.188983712426
$wsCumpl.Range("F7").Value = 2.650130602683176

# PIEDRA SINTERIZADA (row 11)
$wsCumpl.Range("D11").Value = 12981.34
$wsCumpl.Range("E11").Value = 4850.0743984654
$wsCumpl.Range("F11").Value = 0.728003943485111

# PORCELANATO (row 12)
$wsCumpl.Range("D12").Value = 29357.76
$wsCumpl.Range("E12").Value = 32505.9603947566
$wsCumpl.Range("F12").Value = 0.4745553583371019

# SAL SOLUBLE (row 14)
$wsCumpl.Range("D14").Value = 3643.45
$wsCumpl.Range("E14").Value = 4193.86410570622
$wsCumpl.Range("F14").Value = 0.4648850295979924

# TOTAL (row 15)
$wsCumpl.Range("D15").Value = 57226.82999999999
$wsCumpl.Range("E15").Value = 64828.00551083435
$wsCumpl.Range("F15").Value = 0.468861637152591
